$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "CRMs shock data"
$ws.Range("B11").Value = "C:\Users\loren\Documents\GitHub\SESAM\GreenTechs\Support data\CRMs shocks.xlsx"

$ws.Range("A12").Select()
